$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo: "Seeed Weight (grams)" -> "Seed Weight (grams)"
$ws.Range("B8").Value = "Seed Weight (grams)"

# Reflect the last user selection being on B8 (as in the target file)
$ws.Range("B8").Select()
